# Populate the "KIT COMPONENTS" table (3rd table in the document) with the
# actual kit component rows, replacing the placeholder/overview text that
# was previously duplicated into this table.

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(3)

# Row 2 (index 1 data row): Anti-Mouse Klk1 Pre-coated 96-well Strip Microplate
$tbl.Cell(2, 1).Range.Text = "Anti-Mouse Klk1 Pre-coated 96-well Strip Microplate"
$tbl.Cell(2, 2).Range.Text = "1"
$tbl.Cell(2, 3).Range.Text = "12 strips of 8 wells"
$tbl.Cell(2, 4).Range.Text = "Return unused wells to the foil pouch. Reseal along the entire edge of the zip-seal. May be stored for up to 1 month at 4°C provided this is within the expiration date of the kit."

# Row 3: Mouse Klk1 Standard
$tbl.Cell(3, 1).Range.Text = "Mouse Klk1 Standard"
$tbl.Cell(3, 2).Range.Text = "2"
$tbl.Cell(3, 3).Range.Text = "10 ng/tube"
$tbl.Cell(3, 4).Range.Text = "Discard the Klk1 stock solution after 12 hours at 4°C. May be stored at -20°C for 48 hours."

# Row 4: Mouse Klk1 Biotinylated Antibody (100x)
$tbl.Cell(4, 1).Range.Text = "Mouse Klk1 Biotinylated Antibody (100x)"
$tbl.Cell(4, 2).Range.Text = "1"
$tbl.Cell(4, 3).Range.Text = "100 µl"
$tbl.Cell(4, 4).Range.Text = "May be stored for up to 1 month at 4°C provided this is within the expiration date of the kit."

# Row 5: Avidin-Biotin-Peroxidase Complex (100x)
$tbl.Cell(5, 1).Range.Text = "Avidin-Biotin-Peroxidase Complex (100x)"
$tbl.Cell(5, 2).Range.Text = "1"
$tbl.Cell(5, 3).Range.Text = "100 µl"
$tbl.Cell(5, 4).Range.Text = "May be stored for up to 1 month at 4°C provided this is within the expiration date of the kit."

# Row 6: Sample Diluent
$tbl.Cell(6, 1).Range.Text = "Sample Diluent"
$tbl.Cell(6, 2).Range.Text = "1"
$tbl.Cell(6, 3).Range.Text = "30 ml"
$tbl.Cell(6, 4).Range.Text = "May be stored for up to 1 month at 4°C provided this is within the expiration date of the kit."

# Row 7: Antibody Diluent
$tbl.Cell(7, 1).Range.Text = "Antibody Diluent"
$tbl.Cell(7, 2).Range.Text = "1"
$tbl.Cell(7, 3).Range.Text = "12 ml"
$tbl.Cell(7, 4).Range.Text = "May be stored for up to 1 month at 4°C provided this is within the expiration date of the kit."

# Row 8: Avidin-Biotin-Peroxidase Diluent
$tbl.Cell(8, 1).Range.Text = "Avidin-Biotin-Peroxidase Diluent"
$tbl.Cell(8, 2).Range.Text = "1"
$tbl.Cell(8, 3).Range.Text = "12 ml"
$tbl.Cell(8, 4).Range.Text = "May be stored for up to 1 month at 4°C provided this is within the expiration date of the kit."

Write-Host "Kit components table updated."
